# Add two new columns, "I0" (I) and "IF" (J), to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting already used by B1:H1 (bold font, thin
# border, centered/top aligned) by copying the format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data ----------------------------------------------------------------
$values = @{
    2  = 9
    3  = 9
    4  = 9
    5  = 9
    6  = 9
    7  = 9
    8  = 9
    9  = 9
    10 = 9
    11 = 9
    12 = 9
    13 = 9
    14 = 9
    15 = 10
    16 = 9
    17 = 9
    18 = 9
    19 = 9
    20 = 9
    21 = 9
    22 = 9
    23 = 9
    24 = 8
    25 = 9
    26 = 4
    27 = 3
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("I$row").Value = $v
    $ws.Range("J$row").Value = $v
}
